$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the formula-like text cells as explicit text (quote-prefixed),
# so Excel stores them with quotePrefix="1" style instead of trying to
# evaluate them as formulas.
$ws.Range("D4").Value = "'" + '=  CoveragePremiumCalculation (  new String[0] )'
$ws.Range("D5").Value = "'" + '= (Object[]) $Rate$TotalCoveragePremium ( $RaterCoverages )'
$ws.Range("D6").Value = "'" + '= sum ( $Totals )'

# Move the active selection, matching where the user ended up editing.
$null = $ws.Range("C14:D14").Select()
